$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1045.75
$ws.Range("I107").Value = 751
$ws.Range("J107").Value = 1340.5
$ws.Range("K107").Value = 751
$ws.Range("L107").Value = 1340.5
$ws.Range("M107").Value = 1169
$ws.Range("N107").Value = -5180.5
$ws.Range("H112").Value = 3686.879
$ws.Range("J112").Value = 3686.879
$ws.Range("L112").Value = 11060.637
$ws.Range("N112").Value = -13276.637
$ws.Range("H129").Value = 1062.1444
$ws.Range("J129").Value = 1105.9529
$ws.Range("L129").Value = 3317.8587
$ws.Range("N129").Value = -13317.8587
$ws.Range("H137").Value = 3051.8096
$ws.Range("I137").Value = 2097.1428
$ws.Range("J137").Value = 3529.1428
$ws.Range("K137").Value = 6291.428400000001
$ws.Range("L137").Value = 10587.4284
$ws.Range("M137").Value = -3741.428400000001
$ws.Range("N137").Value = -15687.4284
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = 0
$ws.Range("H61").Value = 2751.6453
$ws.Range("I61").Value = 2774.0435
$ws.Range("K61").Value = 2774.0435
$ws.Range("M61").Value = -2562.0435
$ws.Range("H74").Value = 1512.2554
$ws.Range("I74").Value = 1339.3143
$ws.Range("J74").Value = 2016.6666
$ws.Range("K74").Value = 1339.3143
$ws.Range("L74").Value = 2016.6666
$ws.Range("M74").Value = -465.3143
$ws.Range("N74").Value = -3764.6666
$ws.Range("H77").Value = 1512.2554
$ws.Range("I77").Value = 1339.3143
$ws.Range("J77").Value = 2016.6666
$ws.Range("K77").Value = 6696.5715
$ws.Range("L77").Value = 10083.333
$ws.Range("M77").Value = -2328.5715
$ws.Range("N77").Value = -18819.333
$ws.Range("H102").Value = 2545.5557
$ws.Range("I102").Value = 2701.4285
$ws.Range("K102").Value = 2701.4285
$ws.Range("M102").Value = -1079.4285
$ws.Range("H136").Value = 2751.6453
$ws.Range("I136").Value = 2774.0435
$ws.Range("K136").Value = 8322.130500000001
$ws.Range("M136").Value = -5772.130500000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3295.7778
$ws.Range("J134").Value = 3500
$ws.Range("L134").Value = 10500
$ws.Range("N134").Value = -15570
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9051.593000000001
$ws.Range("J4").Value = 9949.708000000001
$ws.Range("L4").Value = 9949.708000000001
$ws.Range("N4").Value = -10173.708
$ws.Range("H31").Value = 1927.5128
$ws.Range("I31").Value = 1329.625
$ws.Range("J31").Value = 2884.1333
$ws.Range("K31").Value = 1329.625
$ws.Range("L31").Value = 2884.1333
$ws.Range("M31").Value = -1034.625
$ws.Range("N31").Value = -3474.1333
$ws.Range("H34").Value = 1927.5128
$ws.Range("I34").Value = 1329.625
$ws.Range("J34").Value = 2884.1333
$ws.Range("K34").Value = 1329.625
$ws.Range("L34").Value = 2884.1333
$ws.Range("M34").Value = -1127.625
$ws.Range("N34").Value = -3288.1333
$ws.Range("H58").Value = 1157.5581
$ws.Range("I58").Value = 1114.4849
$ws.Range("J58").Value = 1299.7
$ws.Range("K58").Value = 1114.4849
$ws.Range("L58").Value = 1299.7
$ws.Range("M58").Value = -911.4848999999999
$ws.Range("N58").Value = -1705.7
$ws.Range("H86").Value = 4045
$ws.Range("I86").Value = 3429.2856
$ws.Range("K86").Value = 3429.2856
$ws.Range("M86").Value = -2306.2856
$ws.Range("H89").Value = 4045
$ws.Range("I89").Value = 3429.2856
$ws.Range("K89").Value = 17146.428
$ws.Range("M89").Value = -11530.428
$ws.Range("H132").Value = 3055.4138
$ws.Range("I132").Value = 2688.28
$ws.Range("K132").Value = 8064.84
$ws.Range("M132").Value = -5534.84
$ws.Range("H134").Value = 2155.5334
$ws.Range("I134").Value = 1757.5454
$ws.Range("J134").Value = 3250
$ws.Range("K134").Value = 5272.6362
$ws.Range("L134").Value = 9750
$ws.Range("M134").Value = -2737.6362
$ws.Range("N134").Value = -14820
$ws.Range("H136").Value = 1157.5581
$ws.Range("I136").Value = 1114.4849
$ws.Range("J136").Value = 1299.7
$ws.Range("K136").Value = 3343.4547
$ws.Range("L136").Value = 3899.1
$ws.Range("M136").Value = -793.4546999999998
$ws.Range("N136").Value = -8999.1
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 65098.562
$ws.Range("I4").Value = 166749.5
$ws.Range("J4").Value = 4108
$ws.Range("K4").Value = 500248.5
$ws.Range("L4").Value = 12324
$ws.Range("M4").Value = -500136.5
$ws.Range("N4").Value = -12548
$ws.Range("H12").Value = 160
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 180
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 540
$ws.Range("M12").Value = -127
$ws.Range("N12").Value = -886
$ws.Range("H68").Value = 173200.9
$ws.Range("I68").Value = 294683.34
$ws.Range("J68").Value = 1100.7916
$ws.Range("K68").Value = 884050.02
$ws.Range("L68").Value = 3302.3748
$ws.Range("M68").Value = -883239.02
$ws.Range("N68").Value = -4924.3748
$ws.Range("H71").Value = 173200.9
$ws.Range("I71").Value = 294683.34
$ws.Range("J71").Value = 1100.7916
$ws.Range("K71").Value = 2652150.06
$ws.Range("L71").Value = 9907.124400000001
$ws.Range("M71").Value = -2648094.06
$ws.Range("N71").Value = -18019.1244
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 28893.334
$ws.Range("J5").Value = 30005
$ws.Range("L5").Value = 30005
$ws.Range("N5").Value = -30229
$ws.Range("H122").Value = 2358.5
$ws.Range("I122").Value = 2364.1667
$ws.Range("J122").Value = 2350
$ws.Range("K122").Value = 7092.500100000001
$ws.Range("L122").Value = 7050
$ws.Range("M122").Value = -4642.500100000001
$ws.Range("N122").Value = -11950
$ws.Range("H126").Value = 2897.2307
$ws.Range("I126").Value = 2512.6667
$ws.Range("K126").Value = 7538.000100000001
$ws.Range("M126").Value = -5068.000100000001
$ws.Range("H132").Value = 2777.76
$ws.Range("I132").Value = 2507.3333
$ws.Range("J132").Value = 3473.1428
$ws.Range("K132").Value = 7521.999899999999
$ws.Range("L132").Value = 10419.4284
$ws.Range("M132").Value = -4991.999899999999
$ws.Range("N132").Value = -15479.4284
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9356
$ws.Range("I2").Value = 8900
$ws.Range("K2").Value = 8900
$ws.Range("M2").Value = -8788
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H132").Value = 3284.347
$ws.Range("I132").Value = 3317.568
$ws.Range("J132").Value = 2992
$ws.Range("K132").Value = 9952.704000000002
$ws.Range("L132").Value = 8976
$ws.Range("M132").Value = -7422.704000000002
$ws.Range("N132").Value = -14036
$ws.Range("H136").Value = 1043.5416
$ws.Range("I136").Value = 733.9474
$ws.Range("J136").Value = 2220
$ws.Range("K136").Value = 2201.8422
$ws.Range("L136").Value = 6660
$ws.Range("M136").Value = 348.1578
$ws.Range("N136").Value = -11760
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4688
$ws.Range("I62").Value = 4666.6665
$ws.Range("J62").Value = 4699.636
$ws.Range("K62").Value = 4666.6665
$ws.Range("L62").Value = 4699.636
$ws.Range("M62").Value = -4042.6665
$ws.Range("N62").Value = -5947.636
$ws.Range("H65").Value = 4688
$ws.Range("I65").Value = 4666.6665
$ws.Range("J65").Value = 4699.636
$ws.Range("K65").Value = 23333.3325
$ws.Range("L65").Value = 23498.18
$ws.Range("M65").Value = -20213.3325
$ws.Range("N65").Value = -29738.18
$ws.Range("H132").Value = 1999.9
$ws.Range("I132").Value = 1630.5873
$ws.Range("J132").Value = 3368.5293
$ws.Range("K132").Value = 4891.7619
$ws.Range("L132").Value = 10105.5879
$ws.Range("M132").Value = -2361.7619
$ws.Range("N132").Value = -15165.5879
$ws.Range("H136").Value = 1647.9459
$ws.Range("I136").Value = 1399.1923
$ws.Range("J136").Value = 2235.9092
$ws.Range("K136").Value = 4197.5769
$ws.Range("L136").Value = 6707.7276
$ws.Range("M136").Value = -1647.5769
$ws.Range("N136").Value = -11807.7276

Write-Host "Applied all cell updates"